$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 280
$ws.Range("I6").Value = 280
$ws.Range("K6").Value = 840
$ws.Range("M6").Value = -728
$ws.Range("H13").Value = 683.3333
$ws.Range("I13").Value = 525
$ws.Range("K13").Value = 525
$ws.Range("M13").Value = -356
$ws.Range("H15").Value = 366.60416
$ws.Range("I15").Value = 366.60416
$ws.Range("K15").Value = 1099.81248
$ws.Range("M15").Value = -930.8124800000001
$ws.Range("H38").Value = 443.66666
$ws.Range("I38").Value = 443.66666
$ws.Range("K38").Value = 1330.99998
$ws.Range("M38").Value = -958.9999800000001
$ws.Range("H42").Value = 112.22222
$ws.Range("I42").Value = 112.22222
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 336.66666
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -106.66666
$ws.Range("N42").ClearContents()
$ws.Range("H58").Value = 4248.75
$ws.Range("I58").Value = 497.5
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 1492.5
$ws.Range("L58").Value = 24000
$ws.Range("M58").Value = -1342.5
$ws.Range("N58").Value = -24300
$ws.Range("H64").Value = 5887.5
$ws.Range("I64").Value = 7500
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 7500
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -7252
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 5887.5
$ws.Range("I67").Value = 7500
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 7500
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -6642
$ws.Range("N67").Value = -4916
$ws.Range("H82").Value = 324.75
$ws.Range("I82").Value = 324.75
$ws.Range("K82").Value = 974.25
$ws.Range("M82").Value = -568.25
$ws.Range("H85").Value = 324.75
$ws.Range("I85").Value = 324.75
$ws.Range("K85").Value = 974.25
$ws.Range("M85").Value = 429.75
$ws.Range("H92").Value = 212
$ws.Range("I92").Value = 199
$ws.Range("K92").Value = 199
$ws.Range("M92").Value = 1049
$ws.Range("H94").Value = 4780.8
$ws.Range("I94").Value = 4780.8
$ws.Range("K94").Value = 4780.8
$ws.Range("M94").Value = -4329.8
$ws.Range("H98").Value = 1023.4286
$ws.Range("I98").Value = 1024.6
$ws.Range("K98").Value = 1024.6
$ws.Range("M98").Value = 473.4000000000001
$ws.Range("H122").Value = 1023.4286
$ws.Range("I122").Value = 1024.6
$ws.Range("K122").Value = 3073.8
$ws.Range("M122").Value = -623.7999999999997
$ws.Range("H129").Value = 909.6667
$ws.Range("I129").Value = 606.875
$ws.Range("J129").Value = 3332
$ws.Range("K129").Value = 1820.625
$ws.Range("L129").Value = 9996
$ws.Range("M129").Value = 3179.375
$ws.Range("N129").Value = -19996
$ws.Range("H135").Value = 1435.05
$ws.Range("I135").Value = 965.0714
$ws.Range("K135").Value = 8685.642600000001
$ws.Range("M135").Value = -6150.642600000001
$ws.Range("H137").Value = 1841.1538
$ws.Range("I137").Value = 1643.25
$ws.Range("K137").Value = 4929.75
$ws.Range("M137").Value = -2379.75
$ws.Range("H138").Value = 2350.2444
$ws.Range("I138").Value = 1857.64
$ws.Range("K138").Value = 5572.92
$ws.Range("M138").Value = -432.9200000000001
$ws.Range("H141").Value = 3644.1853
$ws.Range("I141").Value = 1674.4348
$ws.Range("K141").Value = 5023.3044
$ws.Range("M141").Value = 156.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 4700
$ws.Range("I14").Value = 2500
$ws.Range("J14").Value = 6900
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 6900
$ws.Range("N14").Value = -7250
$ws.Range("M14").Value = -2325
$ws.Range("H16").Value = 4332.3335
$ws.Range("I16").Value = 4332.3335
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4332.3335
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -4045.3335
$ws.Range("N16").ClearContents()
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20616
$ws.Range("H30").Value = 969.375
$ws.Range("I30").Value = 626
$ws.Range("J30").Value = 1999.5
$ws.Range("K30").Value = 626
$ws.Range("L30").Value = 1999.5
$ws.Range("M30").Value = -476
$ws.Range("N30").Value = -2299.5
$ws.Range("H32").Value = 5065.803
$ws.Range("I32").Value = 2583.5386
$ws.Range("K32").Value = 2583.5386
$ws.Range("M32").Value = -2296.5386
$ws.Range("H41").Value = 2601.6667
$ws.Range("I41").Value = 2601.6667
$ws.Range("K41").Value = 2601.6667
$ws.Range("M41").Value = -2187.6667
$ws.Range("H61").Value = 2130.25
$ws.Range("I61").Value = 2148.2104
$ws.Range("J61").Value = 1789
$ws.Range("K61").Value = 2148.2104
$ws.Range("L61").Value = 1789
$ws.Range("M61").Value = -1936.2104
$ws.Range("N61").Value = -2213
$ws.Range("H74").Value = 651.7778
$ws.Range("I74").Value = 543.05884
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 543.05884
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = 330.94116
$ws.Range("N74").Value = -4248
$ws.Range("H76").Value = 58666.332
$ws.Range("J76").Value = 58666.332
$ws.Range("L76").Value = 58666.332
$ws.Range("N76").Value = -59342.332
$ws.Range("H77").Value = 651.7778
$ws.Range("I77").Value = 543.05884
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 2715.2942
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = 1652.7058
$ws.Range("N77").Value = -21236
$ws.Range("H79").Value = 58666.332
$ws.Range("J79").Value = 58666.332
$ws.Range("L79").Value = 58666.332
$ws.Range("N79").Value = -61006.332
$ws.Range("H122").Value = 9499.714
$ws.Range("I122").Value = 1849.75
$ws.Range("J122").Value = 19699.666
$ws.Range("K122").Value = 5549.25
$ws.Range("L122").Value = 59098.99800000001
$ws.Range("M122").Value = -3099.25
$ws.Range("N122").Value = -63998.99800000001
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("M132").Value = -12470
$ws.Range("H136").Value = 2130.25
$ws.Range("I136").Value = 2148.2104
$ws.Range("J136").Value = 1789
$ws.Range("K136").Value = 6444.6312
$ws.Range("L136").Value = 5367
$ws.Range("M136").Value = -3894.6312
$ws.Range("N136").Value = -10467

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 600
$ws.Range("I7").Value = 466.66666
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 466.66666
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -353.66666
$ws.Range("N7").Value = -1226
$ws.Range("H86").Value = 3204.4285
$ws.Range("I86").Value = 3238.5
$ws.Range("K86").Value = 3238.5
$ws.Range("M86").Value = -2115.5
$ws.Range("H89").Value = 3204.4285
$ws.Range("I89").Value = 3238.5
$ws.Range("K89").Value = 16192.5
$ws.Range("M89").Value = -10576.5
$ws.Range("H97").Value = 15222.5
$ws.Range("I97").Value = 15222.5
$ws.Range("K97").Value = 15222.5
$ws.Range("M97").Value = -14231.5
$ws.Range("H105").Value = 5558566.5
$ws.Range("I105").Value = 11907100
$ws.Range("K105").Value = 11907100
$ws.Range("M105").Value = -11905353
$ws.Range("H134").Value = 2298.3
$ws.Range("I134").Value = 1623
$ws.Range("K134").Value = 4869
$ws.Range("M134").Value = -2334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2060
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 2120
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 2120
$ws.Range("N29").Value = -2706
$ws.Range("M29").Value = -1707
$ws.Range("H45").Value = 16642
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 16642
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 16642
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -17828
$ws.Range("H47").Value = 18000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H58").Value = 2598.5652
$ws.Range("I58").Value = 1362.4375
$ws.Range("K58").Value = 1362.4375
$ws.Range("M58").Value = -1159.4375
$ws.Range("H99").Value = 2767.6667
$ws.Range("I99").Value = 1946.6666
$ws.Range("J99").Value = 4409.6665
$ws.Range("K99").Value = 1946.6666
$ws.Range("L99").Value = 4409.6665
$ws.Range("M99").Value = -448.6666
$ws.Range("N99").Value = -7405.6665
$ws.Range("H105").Value = 2587.0833
$ws.Range("I105").Value = 1190.091
$ws.Range("K105").Value = 1190.091
$ws.Range("M105").Value = 556.9090000000001
$ws.Range("H107").Value = 33333674
$ws.Range("J107").Value = 460.66666
$ws.Range("L107").Value = 460.66666
$ws.Range("N107").Value = -4300.66666
$ws.Range("H122").Value = 1049.2858
$ws.Range("J122").Value = 1400
$ws.Range("L122").Value = 4200
$ws.Range("N122").Value = -9100
$ws.Range("H126").Value = 2767.6667
$ws.Range("I126").Value = 1946.6666
$ws.Range("J126").Value = 4409.6665
$ws.Range("K126").Value = 5839.9998
$ws.Range("L126").Value = 13228.9995
$ws.Range("M126").Value = -3369.9998
$ws.Range("N126").Value = -18168.9995
$ws.Range("H132").Value = 1632.0435
$ws.Range("I132").Value = 1129.8889
$ws.Range("J132").Value = 3439.8
$ws.Range("K132").Value = 3389.6667
$ws.Range("L132").Value = 10319.4
$ws.Range("M132").Value = -859.6666999999998
$ws.Range("N132").Value = -15379.4
$ws.Range("H136").Value = 2598.5652
$ws.Range("I136").Value = 1362.4375
$ws.Range("K136").Value = 4087.3125
$ws.Range("M136").Value = -1537.3125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 478.75
$ws.Range("I18").Value = 404.2857
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1212.8571
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -1043.8571
$ws.Range("N18").Value = -3338
$ws.Range("H29").Value = 338517.16
$ws.Range("I29").Value = 1000007.5
$ws.Range("J29").Value = 7772
$ws.Range("K29").Value = 3000022.5
$ws.Range("L29").Value = 23316
$ws.Range("M29").Value = -2999745.5
$ws.Range("N29").Value = -23870
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H80").Value = 2998.5
$ws.Range("J80").Value = 2998.5
$ws.Range("L80").Value = 8995.5
$ws.Range("N80").Value = -10867.5
$ws.Range("H83").Value = 2998.5
$ws.Range("J83").Value = 2998.5
$ws.Range("L83").Value = 26986.5
$ws.Range("N83").Value = -36346.5
$ws.Range("H109").Value = 1166.6666
$ws.Range("I109").Value = 1166.6666
$ws.Range("K109").Value = 3499.9998
$ws.Range("M109").Value = -2459.9998
$ws.Range("H118").Value = 2222
$ws.Range("I118").Value = 2222
$ws.Range("K118").Value = 6666
$ws.Range("M118").Value = -5423
$ws.Range("H139").Value = 1599.6666
$ws.Range("I139").Value = 1599.6666
$ws.Range("K139").Value = 4798.9998
$ws.Range("M139").Value = 341.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6707.15
$ws.Range("I102").Value = 7009.1177
$ws.Range("J102").Value = 4996
$ws.Range("K102").Value = 7009.1177
$ws.Range("L102").Value = 4996
$ws.Range("M102").Value = -5387.1177
$ws.Range("N102").Value = -8240
$ws.Range("H107").Value = 4625.5
$ws.Range("I107").Value = 3001
$ws.Range("J107").Value = 6250
$ws.Range("K107").Value = 3001
$ws.Range("L107").Value = 6250
$ws.Range("M107").Value = -1081
$ws.Range("N107").Value = -10090
$ws.Range("H122").Value = 93228.82000000001
$ws.Range("I122").Value = 1939.875
$ws.Range("K122").Value = 5819.625
$ws.Range("M122").Value = -3369.625
$ws.Range("H123").Value = 34954.6
$ws.Range("J123").Value = 34968.25
$ws.Range("L123").Value = 34968.25
$ws.Range("N123").Value = -39868.25
$ws.Range("H132").Value = 2400
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H16").Value = 1186.0667
$ws.Range("I16").Value = 1183.3636
$ws.Range("J16").Value = 1193.5
$ws.Range("K16").Value = 1183.3636
$ws.Range("L16").Value = 1193.5
$ws.Range("M16").Value = -1013.3636
$ws.Range("N16").Value = -1533.5
$ws.Range("H40").Value = 2512.2856
$ws.Range("I40").Value = 2514.3333
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 2514.3333
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -2378.3333
$ws.Range("N40").Value = -2772
$ws.Range("H43").Value = 8912.875
$ws.Range("J43").Value = 9715.333000000001
$ws.Range("L43").Value = 9715.333000000001
$ws.Range("N43").Value = -10101.333
$ws.Range("H46").Value = 336332.66
$ws.Range("I46").Value = 5999
$ws.Range("J46").Value = 501499.5
$ws.Range("K46").Value = 5999
$ws.Range("L46").Value = 501499.5
$ws.Range("M46").Value = -5811
$ws.Range("N46").Value = -501875.5
$ws.Range("H68").Value = 3133
$ws.Range("J68").Value = 5003
$ws.Range("L68").Value = 5003
$ws.Range("N68").Value = -6501
$ws.Range("H71").Value = 3133
$ws.Range("J71").Value = 5003
$ws.Range("L71").Value = 25015
$ws.Range("N71").Value = -32503
$ws.Range("H122").Value = 3859.5
$ws.Range("I122").Value = 3750
$ws.Range("K122").Value = 11250
$ws.Range("M122").Value = -8800
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 4475
$ws.Range("I132").Value = 4475
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13425
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10895
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 12500
$ws.Range("J20").Value = 12500
$ws.Range("L20").Value = 12500
$ws.Range("N20").Value = -12980
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H62").Value = 7499.8335
$ws.Range("I62").Value = 6500.5
$ws.Range("J62").Value = 7999.5
$ws.Range("K62").Value = 6500.5
$ws.Range("L62").Value = 7999.5
$ws.Range("M62").Value = -5876.5
$ws.Range("N62").Value = -9247.5
$ws.Range("H65").Value = 7499.8335
$ws.Range("I65").Value = 6500.5
$ws.Range("J65").Value = 7999.5
$ws.Range("K65").Value = 32502.5
$ws.Range("L65").Value = 39997.5
$ws.Range("M65").Value = -29382.5
$ws.Range("N65").Value = -46237.5
$ws.Range("I81").Value = 4491.3335
$ws.Range("J81").Value = 2999.5
$ws.Range("K81").Value = 8982.666999999999
$ws.Range("L81").Value = 5999
$ws.Range("M81").Value = -7921.666999999999
$ws.Range("N81").Value = -8121
$ws.Range("I84").Value = 4491.3335
$ws.Range("J84").Value = 2999.5
$ws.Range("K84").Value = 44913.335
$ws.Range("L84").Value = 29995
$ws.Range("M84").Value = -39609.335
$ws.Range("N84").Value = -40603
$ws.Range("H100").Value = 3258.2856
$ws.Range("I100").Value = 2960.6
$ws.Range("K100").Value = 5921.2
$ws.Range("M100").Value = -5380.2
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H122").Value = 1884.4375
$ws.Range("I122").Value = 1904
$ws.Range("K122").Value = 5712
$ws.Range("M122").Value = -3262
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960
$ws.Range("H129").Value = 81922.5
$ws.Range("I129").Value = 81922.5
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 81922.5
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -76922.5
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 2605.5264
$ws.Range("I132").Value = 1867.2
$ws.Range("K132").Value = 5601.6
$ws.Range("M132").Value = -3071.6
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1414.1052
$ws.Range("I136").Value = 941.3333
$ws.Range("K136").Value = 2823.9999
$ws.Range("M136").Value = -273.9998999999998
